$d = $word.ActiveDocument

# --- Replace stub placeholder values with the actual report values ---

# "Отчет по лабораторной работе № NaN" -> "...№8"
$d.Content.Find.Execute(" NaN", $true, $false, $false, $false, $false, $true, 1, $false, "8", 2)

# "ТЕМА" placeholder -> real thesis title
$d.Content.Find.Execute("ТЕМА", $true, $false, $false, $false, $false, $true, 1, $false, "Создание визуального интерфейса для базы данных", 2)

# --- Add the "Table Contents" / "Table Heading" paragraph styles ---

$tableContents = $d.Styles.Add("Table Contents", 1)
$tableContents.BaseStyle = $d.Styles("Normal")
$tableContents.QuickStyle = $true
$tableContents.ParagraphFormat.WidowControl = $false
$tableContents.ParagraphFormat.NoLineNumber = $true

$tableHeading = $d.Styles.Add("Table Heading", 1)
$tableHeading.BaseStyle = $d.Styles("TableContents")
$tableHeading.QuickStyle = $true
$tableHeading.ParagraphFormat.NoLineNumber = $true
$tableHeading.ParagraphFormat.Alignment = 1
$tableHeading.Font.Bold = $true
$tableHeading.Font.BoldBi = $true
